$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh: update price (D) and 1h volume change (E) columns,
# plus a few coin name/link/price swaps where ranking order changed (B/C/D/E).

# Row 2
$ws.Range("D2").Value = '66.924.15'
$ws.Range("E2").Value = '  -1.99%  '

# Row 3
$ws.Range("D3").Value = '2.466.14'
$ws.Range("E3").Value = '  -3.28%  '

# Row 4
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.10'
$ws.Range("E5").Value = '  -1.99%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.23'
$ws.Range("E6").Value = '  -3.83%  '

# Row 7
$ws.Range("E7").Value = '  +0.12%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.511'
$ws.Range("E8").Value = '  -2.85%  '

# Row 9
$ws.Range("D9").Value = '2.465.84'
$ws.Range("E9").Value = '  -3.26%  '

# Row 10
$ws.Range("E10").Value = '  -3.36%  '

# Row 11
$ws.Range("E11").Value = '  -1.01%  '

# Row 12
$ws.Range("E12").Value = '  -2.83%  '

# Row 13
$ws.Range("E13").Value = '  -5.20%  '

# Row 14
$ws.Range("D14").Value = '2.916.25'
$ws.Range("E14").Value = '  -1.72%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.29'
$ws.Range("E15").Value = '  -4.82%  '

# Row 16
$ws.Range("D16").Value = '66.767.21'
$ws.Range("E16").Value = '  -1.78%  '

# Row 17
$ws.Range("E17").Value = '  -5.11%  '

# Row 18
$ws.Range("D18").Value = '2.484.51'
$ws.Range("E18").Value = '  -2.44%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.96'
$ws.Range("E19").Value = '  -8.02%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.45'
$ws.Range("E20").Value = '  -7.70%  '

# Row 21
$ws.Range("E21").Value = '  -5.84%  '

# Row 22
$ws.Range("E22").Value = '  -3.80%  '

# Row 23
$ws.Range("E23").Value = '  +0.06%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.72'
$ws.Range("E24").Value = '  -4.99%  '

# Row 25
$ws.Range("E25").Value = '  -8.16%  '

# Row 26
$ws.Range("E26").Value = '  -5.77%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.11'
$ws.Range("E27").Value = '  -8.83%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("E28").Value = '  -58.09%  '

# Row 29
$ws.Range("E29").Value = '  -2.57%  '

# Row 30
$ws.Range("D30").Value = '0.0₃0896'
$ws.Range("E30").Value = '  -7.71%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '505.96'
$ws.Range("E31").Value = '  -6.58%  '

# Row 32
$ws.Range("E32").Value = '  -8.28%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.75'
$ws.Range("E33").Value = '  -6.39%  '

# Row 34
$ws.Range("E34").Value = '  -7.78%  '

# Row 35
$ws.Range("E35").Value = '  -0.03%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '159.73'
$ws.Range("E36").Value = '  -0.31%  '

# Row 37
$ws.Range("E37").Value = '  -12.86%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.64'
$ws.Range("E38").Value = '  +0.07%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.20'
$ws.Range("E39").Value = '  -5.79%  '

# Row 40
$ws.Range("E40").Value = '  -9.03%  '

# Row 41
$ws.Range("E41").Value = '  +0.19%  '

# Row 42
$ws.Range("E42").Value = '  -6.23%  '

# Row 43
$ws.Range("B43").Value = 'PolygonEcosystemToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.325'
$ws.Range("E43").Value = '  -7.31%  '

# Row 44
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.77'
$ws.Range("E44").Value = '  -7.58%  '

# Row 45
$ws.Range("E45").Value = '  -7.63%  '

# Row 46
$ws.Range("E46").Value = '  -1.25%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '140.51'
$ws.Range("E47").Value = '  -5.72%  '

# Row 48
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.509'
$ws.Range("E48").Value = '  -8.16%  '

# Row 49
$ws.Range("B49").Value = 'Filecoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.41'
$ws.Range("E49").Value = '  -8.29%  '

# Row 50
$ws.Range("D50").Value = '0.0₆0251'
$ws.Range("E50").Value = '  -11.46%  '

# Row 51
$ws.Range("E51").Value = '  -8.93%  '
